$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "26.413.34"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.610.98"
Set-TextValue "D5" "212.27"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  -0.11%  "
Set-TextValue "D10" "19.22"
$ws.Range("E10").Value = "  +1.27%  "
Set-TextValue "D11" "0.0847"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.837.66"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "1.612.66"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  -0.15%  "
Set-TextValue "D16" "63.57"
Set-TextValue "D17" "234.05"
$ws.Range("E17").Value = "  +8.82%  "
$ws.Range("D18").Value = "26.409.77"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +0.44%  "
Set-TextValue "D20" "7.64"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("E24").Value = "  +0.06%  "
Set-TextValue "D25" "146.82"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +0.23%  "
Set-TextValue "D29" "15.47"
$ws.Range("E29").Value = "  +2.52%  "
Set-TextValue "D30" "0.0496"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "1.505.32"
$ws.Range("E32").Value = "  +6.29%  "
Set-TextValue "D33" "3.24"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("E36").Value = "  -0.22%  "
Set-TextValue "D37" "0.561"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("E38").Value = "  -0.12%  "
Set-TextValue "D39" "0.826"
$ws.Range("E39").Value = "  +0.49%  "
Set-TextValue "D40" "5.85"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("E41").Value = "  -0.02%  "
Set-TextValue "D42" "2.18"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "1.749.88"
$ws.Range("E43").Value = "  +1.20%  "
Set-TextValue "D44" "0.762"
$ws.Range("E44").Value = "  +0.10%  "
Set-TextValue "D45" "0.917"
$ws.Range("E45").Value = "  -1.61%  "
Set-TextValue "D46" "61.35"
$ws.Range("E46").Value = "  +0.76%  "
Set-TextValue "D47" "89.60"
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.49"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.0501"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.0958"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "7.48"
$ws.Range("E51").Value = "  +1.03%  "
